# Daily attendance processing - 2025-12-19 09:03:52
# Reorders the "Recorded By" audit-trail text in column G: the last two
# comma-separated entries in each affected row are swapped (e.g. the most
# recently-appended name moves ahead of the one before it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Value()

    $parts = $text -split ',\s*'
    if ($parts.Count -ge 2) {
        $last = $parts[$parts.Count - 1]
        $secondLast = $parts[$parts.Count - 2]
        $parts[$parts.Count - 1] = $secondLast
        $parts[$parts.Count - 2] = $last
    }

    $cell.Value = [string]::Join(", ", $parts)
}
